$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sports recommendation description to reflect 2 sports instead of 1
$ws.Range("B3").Value = "get a recommendation of sports, return a list of sports (2 now)"

# Update the selection/active cell to B3, matching the recorded view state
$ws.Range("B3").Select()
